$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = "2024-06-15 12:22:20"
$ws.Range("D30").Value = 200
$ws.Range("E30").Value = 9

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = "2024-06-15 12:22:20"
$ws.Range("D31").Value = 200
$ws.Range("E31").Value = 0
